$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.966.59'
$ws.Range("E2").Value = '  -4.82%  '
$ws.Range("D3").Value = '2.206.97'
$ws.Range("E3").Value = '  -7.65%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '294.88'
$ws.Range("E5").Value = '  -6.12%  '
$ws.Range("D6").Value = '79.84'
$ws.Range("E6").Value = '  -10.03%  '
$ws.Range("D7").Value = '0.503'
$ws.Range("E7").Value = '  -5.61%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -8.53%  '
$ws.Range("E10").Value = '  -8.73%  '
$ws.Range("D11").Value = '27.80'
$ws.Range("D12").Value = '45.96'
$ws.Range("E12").Value = '  -12.97%  '
$ws.Range("E13").Value = '  -1.97%  '
$ws.Range("D14").Value = '2.548.86'
$ws.Range("E14").Value = '  -7.61%  '
$ws.Range("E15").Value = '  -8.44%  '
$ws.Range("D16").Value = '13.79'
$ws.Range("E16").Value = '  -9.86%  '
$ws.Range("D17").Value = '2.219.96'
$ws.Range("E17").Value = '  -6.84%  '
$ws.Range("D18").Value = '0.705'
$ws.Range("E18").Value = '  -8.58%  '
$ws.Range("D19").Value = '38.855.21'
$ws.Range("E19").Value = '  -4.95%  '
$ws.Range("D20").Value = '0.0₃0854'
$ws.Range("E20").Value = '  -7.01%  '
$ws.Range("D21").Value = '5.64'
$ws.Range("E21").Value = '  -8.81%  '
$ws.Range("D22").Value = '64.43'
$ws.Range("E22").Value = '  -7.02%  '
$ws.Range("D23").Value = '9.75'
$ws.Range("E23").Value = '  -11.46%  '
$ws.Range("D24").Value = '223.14'
$ws.Range("E24").Value = '  -4.56%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").Value = '2.37'
$ws.Range("E26").Value = '  -10.83%  '
$ws.Range("E27").Value = '  -5.16%  '
$ws.Range("E28").Value = '  -7.75%  '
$ws.Range("E29").Value = '  -2.71%  '
$ws.Range("D30").Value = '8.83'
$ws.Range("E30").Value = '  -6.06%  '
$ws.Range("D31").Value = '147.70'
$ws.Range("E31").Value = '  -3.91%  '
$ws.Range("D32").Value = '30.93'
$ws.Range("E32").Value = '  -9.51%  '
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("E34").Value = '  -10.26%  '
$ws.Range("E35").Value = '  -4.75%  '
$ws.Range("D36").Value = '0.0681'
$ws.Range("E36").Value = '  -7.74%  '
$ws.Range("E37").Value = '  -5.03%  '
$ws.Range("E38").Value = '  -4.35%  '
$ws.Range("D39").Value = '2.60'
$ws.Range("E39").Value = '  -7.60%  '
$ws.Range("E40").Value = '  -9.30%  '
$ws.Range("D41").Value = '14.20'
$ws.Range("E41").Value = '  -12.00%  '
$ws.Range("D42").Value = '3.56'
$ws.Range("E42").Value = '  -7.82%  '
$ws.Range("D43").Value = '1.888.96'
$ws.Range("E43").Value = '  -4.52%  '
$ws.Range("D44").Value = '2.05'
$ws.Range("E44").Value = '  -13.02%  '
$ws.Range("D45").Value = '0.0252'
$ws.Range("E45").Value = '  -7.49%  '
$ws.Range("D46").Value = '15.99'
$ws.Range("E46").Value = '  -9.61%  '
$ws.Range("E47").Value = '  -9.28%  '
$ws.Range("E48").Value = '  -9.97%  '
$ws.Range("D49").Value = '2.423.27'
$ws.Range("E49").Value = '  -7.51%  '
$ws.Range("D50").Value = '69.77'
$ws.Range("E50").Value = '  -4.39%  '
$ws.Range("D51").Value = '1.07'
$ws.Range("E51").Value = '  -1.35%  '
